$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 32, pushing the existing "jpeg_test" row (old row 32)
# down to row 33. Excel's row-insert naturally inherits the formatting of
# the row above (row 31) for the new row 32, and keeps row 33's formatting
# identical to the original row 32 - matching the target workbook exactly.
$ws.Rows("32:32").Insert()

# New row 32: LabVIEW Cyclic Voltammetry entry
$ws.Range("A32").Value = "20.07.2023"
$ws.Range("B32").Value = "CV_50mvs_001"
$ws.Range("C32").Value = "cv_lv"
$ws.Range("I32").Value = "Cyclic Voltammetry ECSTM (LabVIEW)"

# New row 34: LabVIEW Chronoamperometry entry
$ws.Range("A34").Value = "27.02.2018"
$ws.Range("B34").Value = "CA_Pulse-Time-4_01"
$ws.Range("C34").Value = "ca_lv"
$ws.Range("I34").Value = "Chronoamperometry ECSTM (LabVIEW)"

# New row 35: LabVIEW FFT entry (ID/type entered before the date, matching
# the original author's shared-string allocation order)
$ws.Range("B35").Value = "FFT_test"
$ws.Range("C35").Value = "fft_lv"
$ws.Range("A35").Value = "25.05.2023"
$ws.Range("I35").Value = "FFT ECSTM (LabVIEW) - try to show only up to 2 kHz to visualize most relevant noise"

# Re-touch the formerly-blank row 36 placeholder cell: after the writes
# above the engine can drop an all-blank styled row, so nudge its format
# to make sure the wrapped, empty B36 cell (and its row) still round-trips.
$ws.Range("B36").WrapText = $true

# Move the active selection to I38, matching where the editor ended up.
$ws.Range("I38").Select()

Write-Host "edit complete"
